# metadata pipeline logs and staging updates
# Append the new "Austin" / pooled DHM study row to the study_locations sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 5

# Assign new shared strings in the same first-use order as the target
# workbook (StudyID, City, Institution, Study Name) so sharedStrings.xml
# indices line up.
$ws.Cells.Item($row, 1).Value = "DHM Pooled"
$ws.Cells.Item($row, 4).Value = "Austin"
$ws.Cells.Item($row, 8).Value = "Mothers Milk Bank at Austin"
$ws.Cells.Item($row, 3).Value = "Pooled DHM HMO Analysis"

# Match the existing "Analyzed" date column's cell style exactly (built-in
# short-date numFmt) by copying formats from the row above rather than
# assigning a NumberFormat string, which would mint a redundant custom
# numFmt entry in styles.xml.
$ws.Cells.Item(4, 2).Copy()
$ws.Cells.Item($row, 2).PasteSpecial(-4122)
$ws.Cells.Item($row, 2).Value = 45901
$ws.Cells.Item($row, 5).Value = "USA"
$ws.Cells.Item($row, 6).Value = 30.2672
$ws.Cells.Item($row, 7).Value = -97.7431

$ws.Range("G6").Select()
